# Refresh the crypto price/symbol snapshot (GitHub Actions data pull).
# Price cells (column D) are stored as TEXT in this sheet, so every
# numeric-looking value is written with a leading apostrophe to force
# text entry, then ClearFormats() strips the transient "quote prefix"
# number format Excel applies, restoring the cell to its original
# (unstyled) look while keeping the value as text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PriceText {
    param($addr, $text)
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($addr).ClearFormats()
}

# Rows 2-9: price-only refresh
Set-PriceText "D2"  "242.94"
Set-PriceText "D3"  "23.07"
Set-PriceText "D4"  "5.424"
Set-PriceText "D5"  "0.05911"
Set-PriceText "D6"  "3.441"
Set-PriceText "D7"  "6.528"
Set-PriceText "D8"  "0.8097"
Set-PriceText "D9"  "0.9319"

# Rows 10-18: coin listing shifted up one slot (WazirX moves from rank 10
# to rank 9, etc.), each row's Coin/Link/Price/Volume cells updated to the
# next coin in the ranking, wrapping "One" back in at rank 17.
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-PriceText "D10" "0.1427"
$ws.Range("E10").Value = "9WazirXWRX"

$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-PriceText "D11" "0.07414"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"

$ws.Range("B12").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C12").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-PriceText "D12" "0.03236"
$ws.Range("E12").Value = "11LiechtensteinCryptoassetsExchangeLCX"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-PriceText "D13" "0.03093"
$ws.Range("E13").Value = "12BitrueCoinBTR"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-PriceText "D14" "0.09356"
$ws.Range("E14").Value = "13BitMartTokenBMX"

$ws.Range("B15").Value = "MCDex"
$ws.Range("C15").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-PriceText "D15" "3.862"
$ws.Range("E15").Value = "14MCDexMCB"

$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-PriceText "D16" "0.001585"
$ws.Range("E16").Value = "15BitForexTokenBF"

$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-PriceText "D17" "0.04695"
$ws.Range("E17").Value = "16CoinExTokenCET"

$ws.Range("B18").Value = "One"
$ws.Range("C18").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-PriceText "D18" "0.0005917"
$ws.Range("E18").Value = "17OneONE"

# Rows 19-50: remaining price refreshes (and a couple of Volume(1h) label
# tweaks where the "Bestin24h"/"Worstin24h" suffix moved to a different row)
Set-PriceText "D19" "0.005968"
Set-PriceText "D20" "0.001255"
Set-PriceText "D21" "0.004902"
Set-PriceText "D22" "0.00006809"
Set-PriceText "D23" "3.568"
Set-PriceText "D24" "2.140"
Set-PriceText "D26" "0.1304"
Set-PriceText "D27" "0.0002304"
Set-PriceText "D40" "0.03953"

Set-PriceText "D41" "0.003134"
$ws.Range("E41").Value = "40KickTokenKICKWorstin24h"

Set-PriceText "D42" "0.1076"
Set-PriceText "D43" "0.002573"
Set-PriceText "D44" "0.008787"
Set-PriceText "D45" "0.00005214"

Set-PriceText "D47" "0.6708"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"

Set-PriceText "D48" "0.002393"
Set-PriceText "D49" "0.00002103"
Set-PriceText "D50" "0.0002003"
